# Update task data cells used in testing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values: x_corrSteps (D5), y_nrSteps (F5), praclen (H5)
$ws.Range("D5").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("H5").Value = 46

# Move the active selection to D5 (matches the saved selection in the file)
$ws.Range("D5").Select()
